$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "choices" sheet: add a new row for the "families" form, mirroring the
#    existing "marriages" row (both use the "forms" choice_list_name).
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Range("A7").Value = "forms"
$choices.Range("B7").Value = "families"
$choices.Range("C7").Value = "Families"
$choices.Range("A8").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. "survey" sheet: duplicate the "marriages" branch block (rows 30-32) into
#    a new "families" branch block (rows 34-36), keeping formatting/styles.
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("A30:H32").Copy() | Out-Null
$survey.Range("A34:H36").PasteSpecial(-4122) | Out-Null
$survey.Range("A34:H36").PasteSpecial(-4163) | Out-Null

$survey.Range("A34").Value = "families"
$survey.Range("B35").Value = "''?' + odkSurvey.getHashString('families')"

$survey.Range("B36").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add a new, empty "Sheet1" worksheet at the end of the workbook. This
#    becomes the active sheet/tab (mirrors the author having just inserted a
#    blank scratch sheet).
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet1"
